# Auto-generated edit script: refresh market-price derived cells
# per scheduled runner update (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 893085.8
$ws.Range("J17").Value = 893085.8
$ws.Range("L17").Value = 2679257.4
$ws.Range("N17").Value = -2679593.4
# Row 80
$ws.Range("H80").Value = 4546251
$ws.Range("I80").Value = 308.2
$ws.Range("K80").Value = 924.5999999999999
$ws.Range("M80").Value = 73.40000000000009
# Row 83
$ws.Range("H83").Value = 4546251
$ws.Range("I83").Value = 308.2
$ws.Range("K83").Value = 2773.8
$ws.Range("M83").Value = 2218.2
# Row 98
$ws.Range("H98").Value = 1403.1904
$ws.Range("I98").Value = 1360.0555
$ws.Range("J98").Value = 1662
$ws.Range("K98").Value = 1360.0555
$ws.Range("L98").Value = 1662
$ws.Range("M98").Value = 137.9445000000001
$ws.Range("N98").Value = -4658
# Row 107
$ws.Range("H107").Value = 751.3570999999999
$ws.Range("I107").Value = 751.3570999999999
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 751.3570999999999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1168.6429
$ws.Range("N107").ClearContents()
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 122
$ws.Range("H122").Value = 1403.1904
$ws.Range("I122").Value = 1360.0555
$ws.Range("J122").Value = 1662
$ws.Range("K122").Value = 4080.1665
$ws.Range("L122").Value = 4986
$ws.Range("M122").Value = -1630.1665
$ws.Range("N122").Value = -9886
# Row 137
$ws.Range("H137").Value = 3103.5
$ws.Range("I137").Value = 572.5517
$ws.Range("J137").Value = 7996.6665
$ws.Range("K137").Value = 1717.6551
$ws.Range("L137").Value = 23989.9995
$ws.Range("M137").Value = 832.3449000000001
$ws.Range("N137").Value = -29089.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 2088
$ws.Range("I110").Value = 2119.4167
$ws.Range("J110").Value = 1899.5
$ws.Range("K110").Value = 2119.4167
$ws.Range("L110").Value = 1899.5
$ws.Range("M110").Value = -74.41670000000022
$ws.Range("N110").Value = -5989.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 996714.25
$ws.Range("I105").Value = 1991278.5
$ws.Range("J105").Value = 2150
$ws.Range("K105").Value = 1991278.5
$ws.Range("L105").Value = 2150
$ws.Range("M105").Value = -1989531.5
$ws.Range("N105").Value = -5644

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2467.5833
$ws.Range("I62").Value = 2472.1428
$ws.Range("J62").Value = 2461.2
$ws.Range("K62").Value = 2472.1428
$ws.Range("L62").Value = 2461.2
$ws.Range("M62").Value = -1848.1428
$ws.Range("N62").Value = -3709.2
# Row 65
$ws.Range("H65").Value = 2467.5833
$ws.Range("I65").Value = 2472.1428
$ws.Range("J65").Value = 2461.2
$ws.Range("K65").Value = 12360.714
$ws.Range("L65").Value = 12306
$ws.Range("M65").Value = -9240.714
$ws.Range("N65").Value = -18546
# Row 86
$ws.Range("H86").Value = 45464720
$ws.Range("I86").Value = 76936740
$ws.Range("J86").Value = 5133.778
$ws.Range("K86").Value = 76936740
$ws.Range("L86").Value = 5133.778
$ws.Range("M86").Value = -76935617
$ws.Range("N86").Value = -7379.778
# Row 89
$ws.Range("H89").Value = 45464720
$ws.Range("I89").Value = 76936740
$ws.Range("J89").Value = 5133.778
$ws.Range("K89").Value = 384683700
$ws.Range("L89").Value = 25668.89
$ws.Range("M89").Value = -384678084
$ws.Range("N89").Value = -36900.89
# Row 107
$ws.Range("H107").Value = 392.72726
$ws.Range("I107").Value = 392
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 392
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1528
$ws.Range("N107").Value = -4240
# Row 132
$ws.Range("H132").Value = 49398.668
$ws.Range("I132").Value = 67591.734
$ws.Range("J132").Value = 3916
$ws.Range("K132").Value = 202775.202
$ws.Range("L132").Value = 11748
$ws.Range("M132").Value = -200245.202
$ws.Range("N132").Value = -16808

$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 9916.666999999999
$ws.Range("J19").Value = 10500
$ws.Range("L19").Value = 10500
$ws.Range("N19").Value = -11076
# Row 122
$ws.Range("H122").Value = 4275.4287
$ws.Range("I122").Value = 7455.5
$ws.Range("J122").Value = 3003.4
$ws.Range("K122").Value = 22366.5
$ws.Range("L122").Value = 9010.200000000001
$ws.Range("M122").Value = -19916.5
$ws.Range("N122").Value = -13910.2

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1385.1904
$ws.Range("I16").Value = 1361.909
$ws.Range("K16").Value = 1361.909
$ws.Range("M16").Value = -1191.909
# Row 40
$ws.Range("H40").Value = 61409.117
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 74032.5
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 74032.5
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -74304.5
# Row 46
$ws.Range("H46").Value = 1898.5
$ws.Range("I46").Value = 1133.6666
$ws.Range("J46").Value = 2663.3333
$ws.Range("K46").Value = 1133.6666
$ws.Range("L46").Value = 2663.3333
$ws.Range("M46").Value = -945.6666
$ws.Range("N46").Value = -3039.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1920.5
$ws.Range("I81").Value = 2500.5
$ws.Range("J81").Value = 1485.5
$ws.Range("K81").Value = 5001
$ws.Range("L81").Value = 2971
$ws.Range("M81").Value = -3940
$ws.Range("N81").Value = -5093
# Row 84
$ws.Range("H84").Value = 1920.5
$ws.Range("I84").Value = 2500.5
$ws.Range("J84").Value = 1485.5
$ws.Range("K84").Value = 25005
$ws.Range("L84").Value = 14855
$ws.Range("M84").Value = -19701
$ws.Range("N84").Value = -25463
# Row 96
$ws.Range("H96").Value = 3689.6667
$ws.Range("I96").Value = 2201
$ws.Range("J96").Value = 3987.4
$ws.Range("K96").Value = 2201
$ws.Range("L96").Value = 3987.4
$ws.Range("M96").Value = -828
$ws.Range("N96").Value = -6733.4
# Row 107
$ws.Range("H107").Value = 192.6923
$ws.Range("I107").Value = 184
$ws.Range("J107").Value = 212.25
$ws.Range("K107").Value = 552
$ws.Range("L107").Value = 636.75
$ws.Range("M107").Value = 1368
$ws.Range("N107").Value = -4476.75
# Row 113
$ws.Range("H113").Value = 443.7647
$ws.Range("I113").Value = 380.07693
$ws.Range("J113").Value = 650.75
$ws.Range("K113").Value = 1140.23079
$ws.Range("L113").Value = 1952.25
$ws.Range("M113").Value = 1029.76921
$ws.Range("N113").Value = -6292.25
# Row 124
$ws.Range("H124").Value = 16232.25
$ws.Range("J124").Value = 16232.25
$ws.Range("L124").Value = 16232.25
$ws.Range("N124").Value = -26052.25
# Row 126
$ws.Range("H126").Value = 788.8125
$ws.Range("I126").Value = 756.4545000000001
$ws.Range("J126").Value = 860
$ws.Range("K126").Value = 2269.3635
$ws.Range("L126").Value = 2580
$ws.Range("M126").Value = 200.6364999999996
$ws.Range("N126").Value = -7520
